$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "- Taxa de conversão por tipo de evento" -> append " feito" (highlighted)
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Taxa de conversão por tipo de evento", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" feito")

$f = $d.Content
$f.Find.ClearFormatting()
$f.Find.Replacement.ClearFormatting()
$f.Find.Replacement.Highlight = $true
$null = $f.Find.Execute("feito", $true, $false, $false, $false, $false, $true, 1, $true, "feito", 1)

# ---------------------------------------------------------------------------
# 2) "- Produtividade por negociador" -> append " falta colocar a conversão
#    e ver as promessas" (highlighted); the _GoBack bookmark now marks this
#    insertion point.
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("- Produtividade por negociador", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("falta colocar a conversão e ver as promessas")

$f = $d.Content
$f.Find.ClearFormatting()
$f.Find.Replacement.ClearFormatting()
$f.Find.Replacement.Highlight = $true
$null = $f.Find.Execute("falta colocar a conversão e ver as promessas", $true, $false, $false, $false, $false, $true, 1, $true, "falta colocar a conversão e ver as promessas", 1)

# Re-seat the _GoBack bookmark exactly between the space and the new phrase
$b = $d.Content
$null = $b.Find.Execute("negociador ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$b.Collapse(0)
$d.Bookmarks.Add("_GoBack", $b)

# ---------------------------------------------------------------------------
# 3) "- Distribuição por faixa de atraso" -> append " calculado e no gráfico"
#    (highlighted)
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Distribuição por faixa de atraso", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" calculado e no gráfico")

$f = $d.Content
$f.Find.ClearFormatting()
$f.Find.Replacement.ClearFormatting()
$f.Find.Replacement.Highlight = $true
$null = $f.Find.Execute("calculado e no gráfico", $true, $false, $false, $false, $false, $true, 1, $true, "calculado e no gráfico", 1)

Write-Output "done"
